$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide columns B and C, and resize columns A, B, C.
# (ColumnWidth is rounded by the engine to the nearest 1/6 "character" unit
# when serialized back to OOXML <col width=".."/>, so these inputs are the
# values that land closest to the target widths of 16.39 / 23.61 / 28.9.)
$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(2).ColumnWidth = 22.833
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(3).ColumnWidth = 28.0
$ws.Columns.Item(3).Hidden = $false

# Row 2 height
$ws.Rows.Item(2).RowHeight = 23.25

# Move the active selection from N26 to B3
$ws.Range("B3").Select() | Out-Null
